$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Image"
$ws.Range("I3").Value = "아이템 이미지 링크"
$ws.Range("I6").Value = "collectibles_007_bloodofthemartyr"
$ws.Range("I5").Value = "collectibles_006_numberone"
$ws.Range("I4").Value = "collectibles_005_myreflection"
$ws.Range("I7").Value = "collectibles_008_brotherbobby"
$ws.Range("I8").Value = "collectibles_010_haloofflies"
$ws.Range("I9").Value = "collectibles_012_magicmushroom"
$ws.Range("I10").Value = "collectibles_013_thevirus"
$ws.Range("I11").Value = "collectibles_009_skatole"
$ws.Range("I2").Value = "string"

$ws.Rows(3).EntireRow.AutoFit()
